{"js": "// Update each \"A\u00d7B=C\" answer cell to the new value, matched by its\n// current (old) text so the edit is robust to document order.\nconst replacements = [\n  [\"86\u00d767=5762\", \"73\u00d740=2920\"],\n  [\"90\u00d713=1170\", \"70\u00d736=2520\"],\n  [\"19\u00d772=1368\", \"77\u00d723=1771\"],\n  [\"62\u00d730=1860\", \"99\u00d786=8514\"],\n  [\"33\u00d772=2376\", \"52\u00d736=1872\"],\n  [\"71\u00d756=3976\", \"71\u00d728=1988\"],\n  [\"82\u00d798=8036\", \"40\u00d727=1080\"],\n  [\"82\u00d730=2460\", \"53\u00d720=1060\"],\n  [\"66\u00d711=726\", \"78\u00d746=3588\"],\n  [\"22\u00d793=2046\", \"29\u00d729=841\"],\n  [\"49\u00d759=2891\", \"24\u00d785=2040\"],\n  [\"51\u00d733=1683\", \"65\u00d755=3575\"],\n  [\"71\u00d714=994\", \"76\u00d764=4864\"],\n  [\"79\u00d768=5372\", \"47\u00d720=940\"],\n  [\"98\u00d759=5782\", \"13\u00d757=741\"],\n  [\"76\u00d739=2964\", \"33\u00d722=726\"],\n  [\"13\u00d760=780\", \"85\u00d799=8415\"],\n  [\"19\u00d757=1083\", \"62\u00d758=3596\"],\n  [\"47\u00d738=1786\", \"91\u00d769=6279\"],\n  [\"70\u00d740=2800\", \"34\u00d743=1462\"],\n  [\"79\u00d759=4661\", \"87\u00d776=6612\"],\n  [\"62\u00d753=3286\", \"88\u00d748=4224\"],\n  [\"25\u00d786=2150\", \"21\u00d785=1785\"],\n  [\"44\u00d773=3212\", \"96\u00d762=5952\"],\n  [\"87\u00d758=5046\", \"37\u00d747=1739\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text \"${oldText}\" to replace.`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update each \"A\u00d7B=C\" answer cell to the new value, matched by its\n# current (old) text so the edit is robust to document order.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"86\u00d767=5762\", \"73\u00d740=2920\"),\n    @(\"90\u00d713=1170\", \"70\u00d736=2520\"),\n    @(\"19\u00d772=1368\", \"77\u00d723=1771\"),\n    @(\"62\u00d730=1860\", \"99\u00d786=8514\"),\n    @(\"33\u00d772=2376\", \"52\u00d736=1872\"),\n    @(\"71\u00d756=3976\", \"71\u00d728=1988\"),\n    @(\"82\u00d798=8036\", \"40\u00d727=1080\"),\n    @(\"82\u00d730=2460\", \"53\u00d720=1060\"),\n    @(\"66\u00d711=726\",  \"78\u00d746=3588\"),\n    @(\"22\u00d793=2046\", \"29\u00d729=841\"),\n    @(\"49\u00d759=2891\", \"24\u00d785=2040\"),\n    @(\"51\u00d733=1683\", \"65\u00d755=3575\"),\n    @(\"71\u00d714=994\",  \"76\u00d764=4864\"),\n    @(\"79\u00d768=5372\", \"47\u00d720=940\"),\n    @(\"98\u00d759=5782\", \"13\u00d757=741\"),\n    @(\"76\u00d739=2964\", \"33\u00d722=726\"),\n    @(\"13\u00d760=780\",  \"85\u00d799=8415\"),\n    @(\"19\u00d757=1083\", \"62\u00d758=3596\"),\n    @(\"47\u00d738=1786\", \"91\u00d769=6279\"),\n    @(\"70\u00d740=2800\", \"34\u00d743=1462\"),\n    @(\"79\u00d759=4661\", \"87\u00d776=6612\"),\n    @(\"62\u00d753=3286\", \"88\u00d748=4224\"),\n    @(\"25\u00d786=2150\", \"21\u00d785=1785\"),\n    @(\"44\u00d773=3212\", \"96\u00d762=5952\"),\n    @(\"87\u00d758=5046\", \"37\u00d747=1739\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
